# Scoreboard.xlsx update — add the Final round results for FM (men) and
# FF (women), update the selections on SFM/SFF to reflect the qualifying
# teams that advanced, and leave the FF sheet as the active tab/selection
# (matching the state the workbook was saved in).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# FM (Final Men) — append the three teams that made the men's final.
# ---------------------------------------------------------------------
$wsFM = $wb.Worksheets.Item("FM")

$wsFM.Range("A2").Value = "Magnus Øslebye og Vegard Tangen"
$wsFM.Range("B2").Value = 9
$wsFM.Range("C2").Value = 10
$wsFM.Range("D2").Value = 0
$wsFM.Range("E2").Value = 0
$wsFM.Range("F2").Value = 0

$wsFM.Range("A3").Value = "Anders J. Svalestuen og Gabriel Kristiansen"
$wsFM.Range("B3").Value = 9
$wsFM.Range("C3").Value = 12
$wsFM.Range("D3").Value = 0
$wsFM.Range("E3").Value = 0
$wsFM.Range("F3").Value = 0

$wsFM.Range("A4").Value = "Kasper Støen Nerbøvik og Håvard Idland"
$wsFM.Range("B4").Value = 10
$wsFM.Range("C4").Value = 15
$wsFM.Range("D4").Value = 0
$wsFM.Range("E4").Value = 0
$wsFM.Range("F4").Value = 0

# ---------------------------------------------------------------------
# FF (Final Women) — append the three teams that made the women's final.
# ---------------------------------------------------------------------
$wsFF = $wb.Worksheets.Item("FF")

$wsFF.Range("A2").Value = "Dawn Stewart og Marie Vik"
$wsFF.Range("B2").Value = 8
$wsFF.Range("C2").Value = 12
$wsFF.Range("D2").Value = 60
$wsFF.Range("E2").Value = 0
$wsFF.Range("F2").Value = 270

$wsFF.Range("A3").Value = "Sara Yuzer og Martine Baalsrud"
$wsFF.Range("B3").Value = 11
$wsFF.Range("C3").Value = 10
$wsFF.Range("D3").Value = 13
$wsFF.Range("E3").Value = 46
$wsFF.Range("F3").Value = 280

$wsFF.Range("A4").Value = "Frid Kaspersen og Renate Loraas"
$wsFF.Range("B4").Value = 13
$wsFF.Range("C4").Value = 15
$wsFF.Range("D4").Value = 14
$wsFF.Range("E4").Value = 2
$wsFF.Range("F4").Value = 280

# ---------------------------------------------------------------------
# Update the saved cursor/selection on the semi-final sheets (SFM/SFF) to
# highlight the qualifying teams column (A2:A7).
# ---------------------------------------------------------------------
$wsSFM = $wb.Worksheets.Item("SFM")
$wsSFM.Activate()
$wsSFM.Range("A2:A7").Select()

$wsSFF = $wb.Worksheets.Item("SFF")
$wsSFF.Activate()
$wsSFF.Range("A2:A7").Select()

# ---------------------------------------------------------------------
# Leave the selections on FM/FF as saved, and FF as the final active tab.
# ---------------------------------------------------------------------
$wsFM.Activate()
$wsFM.Range("I14").Select()

$wsFF.Activate()
$wsFF.Range("D5").Select()
